$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-05 Saturday" "2025-04-06 Sunday"

Replace-Text "753×8=6024" "882×2=1764"
Replace-Text "783×2=1566" "691×5=3455"
Replace-Text "141×6=846" "411×5=2055"
Replace-Text "424×9=3816" "496×2=992"
Replace-Text "727×6=4362" "542×6=3252"

Replace-Text "216×3=648" "869×6=5214"
Replace-Text "152×3=456" "595×3=1785"
Replace-Text "344×6=2064" "878×6=5268"
Replace-Text "357×4=1428" "421×5=2105"
Replace-Text "306×6=1836" "663×9=5967"

Replace-Text "477×2=954" "465×4=1860"
Replace-Text "331×6=1986" "637×5=3185"
Replace-Text "891×7=6237" "347×4=1388"
Replace-Text "771×4=3084" "274×6=1644"
Replace-Text "631×9=5679" "350×5=1750"

Replace-Text "303×4=1212" "435×3=1305"
Replace-Text "861×2=1722" "390×4=1560"
Replace-Text "800×4=3200" "200×4=800"
Replace-Text "341×2=682" "949×3=2847"
Replace-Text "204×9=1836" "633×4=2532"

Replace-Text "786×7=5502" "745×8=5960"
Replace-Text "264×8=2112" "764×7=5348"
Replace-Text "762×9=6858" "218×6=1308"
Replace-Text "184×7=1288" "519×6=3114"
Replace-Text "339×6=2034" "569×8=4552"
